$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update column C (Förändrad) for existing data rows 2..237 from 45205 -> 45206
$ws.Range("C2:C237").Value = 45206

# Row 237 previously had no explicit row height; set it to match other rows (15pt, custom height)
$ws.Rows.Item(237).RowHeight = 15

# Add new row 238 data
$ws.Range("A238").Value = "A 47832-2023"
$ws.Range("B238").Value = 45204
$ws.Range("C238").Value = 45206
$ws.Range("D238").Value = "HALLANDS LÄN"
$ws.Range("E238").Value = "KUNGSBACKA"
$ws.Range("G238").Value = 2.9
$ws.Range("H238").Value = 0
$ws.Range("I238").Value = 0
$ws.Range("J238").Value = 0
$ws.Range("K238").Value = 0
$ws.Range("L238").Value = 0
$ws.Range("M238").Value = 0
$ws.Range("N238").Value = 0
$ws.Range("O238").Value = 0
$ws.Range("P238").Value = 0
$ws.Range("Q238").Value = 0

# Match styles: B/C use date style (index 1 -> number format YYYY-MM-DD)
$ws.Range("B238:C238").NumberFormat = "YYYY-MM-DD"

# R238 should carry wrap-text style like other rows, left blank
$ws.Range("R238").WrapText = $true
